$wb = $excel.ActiveWorkbook

# --- Sheet "Range Status": zero out column B (species no.) and clear column C (species perc.) ---
$ws = $wb.Worksheets.Item("Range Status")
$ws.Range("B2").Value = 0
$ws.Range("C2").ClearContents()

$ws.Range("C3").ClearContents()

$ws.Range("B4").Value = 0
$ws.Range("C4").ClearContents()

$ws.Range("B5").Value = 0
$ws.Range("C5").ClearContents()

$ws.Range("B6").Value = 0
$ws.Range("C6").ClearContents()

$ws.Range("B7").Value = 0
$ws.Range("C7").ClearContents()

# --- Sheet "Species qualification": Range Analysis selected-for-analysis count resets to 0 ---
$ws = $wb.Worksheets.Item("Species qualification")
$ws.Range("B5").Value = 0

# --- Sheet "High Priority break-up": recompute percentages / add new-high-species breakdown for IUCN ---
$ws = $wb.Worksheets.Item("High Priority break-up")
$ws.Range("E2").Value = 8.300000000000001
$ws.Range("D3").Value = 11
$ws.Range("E3").Value = 91.7
